# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, reflecting the newer scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 2061
$ws.Range("F6").Value  = 4211
$ws.Range("F8").Value  = 1051
$ws.Range("F14").Value = 660119
$ws.Range("F15").Value = 1632
$ws.Range("F18").Value = 671
$ws.Range("F21").Value = 2240
$ws.Range("F25").Value = 818
$ws.Range("F28").Value = 529
$ws.Range("F30").Value = 283
$ws.Range("F37").Value = 1266
$ws.Range("F38").Value = 2535
$ws.Range("F42").Value = 2584
$ws.Range("F45").Value = 3133

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 482
$ws.Range("F11").Value = 144632
$ws.Range("F12").Value = 144632
$ws.Range("F22").Value = 128
$ws.Range("F26").Value = 568
$ws.Range("F37").Value = 3
$ws.Range("F39").Value = 112

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value  = 1182
$ws.Range("F12").Value = 95

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 1182
$ws.Range("F9").Value  = 2061
$ws.Range("F10").Value = 95
$ws.Range("F12").Value = 4211
$ws.Range("F18").Value = 660119
$ws.Range("F20").Value = 482
$ws.Range("F21").Value = 1632
$ws.Range("F22").Value = 144632
$ws.Range("F24").Value = 671
$ws.Range("F27").Value = 2240
$ws.Range("F31").Value = 818
$ws.Range("F34").Value = 529
$ws.Range("F35").Value = 128
$ws.Range("F41").Value = 1266
$ws.Range("F42").Value = 2535
$ws.Range("F48").Value = 2584
$ws.Range("F51").Value = 3133
